$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the effort estimation values in column B
$ws.Range("B4").Value = 50
$ws.Range("B7").Value = 120
$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 30
$ws.Range("B10").Value = 50
$ws.Range("B13").Value = 30
$ws.Range("B14").Value = 30

# Update the active cell selection to D8
$ws.Range("D8").Select()
